# Applies the "Vivalata les modifs !" changes:
#  - adds a 5th column (xb2 / "-" markers) with bound formulas in D/E and helper
#    computations in columns J/K/L
#  - draws a bordered box around the little table (thin box around the header
#    row, medium box around the data rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New header cell + new shared strings ("xb2", "-")
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "xb2"
$ws.Range("E2").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("E7").Value = "-"

# ---------------------------------------------------------------------------
# 2) Helper formulas in columns J, K, L (rows 5-8)
# ---------------------------------------------------------------------------
$ws.Range("J6").Formula = "=256-4*7*(8-3/4*LN(2))"
$ws.Range("K6").Formula = "=2+LN(2)"
$ws.Range("L6").Formula = "=4+4*7*(1+4*LN(2))"

$ws.Range("J5").Formula = "=SQRT(J6)"
$ws.Range("K5").Formula = "=SQRT(K6)"
$ws.Range("L5").Formula = "=SQRT(L6)"

$ws.Range("J7").Formula = "=LN(2)"
$ws.Range("J8").Formula = "=SQRT(3)"

# ---------------------------------------------------------------------------
# 3) D/E formulas for rows 4-6 (bounds around xb / xb2)
# ---------------------------------------------------------------------------
$ws.Range("D4").Formula = "=(32+J5)/2*15"
$ws.Range("E4").Formula = "=(32-J5)/2*15"

$ws.Range("D5").Formula = "=(-1+K5)"
$ws.Range("E5").Formula = "=(-1-K5)"

$ws.Range("D6").Formula = "=(-2+2*L5)/14"
$ws.Range("E6").Formula = "=(-2-2*L5)/14"

$wb.Application.Calculate()

# ---------------------------------------------------------------------------
# 4) Borders - build the "box" around the small table.
#    Processed one cell at a time (same edge-order per cell) so that cells
#    that end up with the same combination of borders share a single style.
# ---------------------------------------------------------------------------
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlThin = 2
$xlMedium = -4138

function Set-Edge($addr, $edge, $weight) {
    $ws.Range($addr).Borders.Item($edge).Weight = $weight
}

# Header row: every cell gets a thin left+right+top border individually.
foreach ($addr in @("A1", "B1", "C1", "D1", "E1")) {
    Set-Edge $addr $xlEdgeLeft $xlThin
    Set-Edge $addr $xlEdgeRight $xlThin
    Set-Edge $addr $xlEdgeTop $xlThin
}

# Row 2 (top of the data box): left+top on A2, top on B2:D2, top+right on E2.
Set-Edge "A2" $xlEdgeLeft $xlMedium
Set-Edge "A2" $xlEdgeTop $xlMedium

foreach ($addr in @("B2", "C2", "D2")) {
    Set-Edge $addr $xlEdgeTop $xlMedium
}

Set-Edge "E2" $xlEdgeTop $xlMedium
Set-Edge "E2" $xlEdgeRight $xlMedium

# Rows 3-6: left border on column A, right border on column E (E6 excluded -
# it keeps no border, matching the source table).
foreach ($r in 3..6) {
    Set-Edge "A$r" $xlEdgeLeft $xlMedium
}
foreach ($r in 3..5) {
    Set-Edge "E$r" $xlEdgeRight $xlMedium
}

# Row 7 (bottom of the data box): left+bottom on A7, bottom on B7:D7,
# bottom+right on E7.
Set-Edge "A7" $xlEdgeLeft $xlMedium
Set-Edge "A7" $xlEdgeBottom $xlMedium

foreach ($addr in @("B7", "C7", "D7")) {
    Set-Edge $addr $xlEdgeBottom $xlMedium
}

Set-Edge "E7" $xlEdgeBottom $xlMedium
Set-Edge "E7" $xlEdgeRight $xlMedium

# Give the empty placeholder cells (C4:C6) the same (border-less) style as
# their neighbours so the shape of the table is consistent.
foreach ($addr in @("C4", "C5", "C6")) {
    $ws.Range($addr).Borders.Item($xlEdgeLeft).LineStyle = -4142
}

# ---------------------------------------------------------------------------
# 5) Misc bookkeeping to mirror the end state of the workbook
# ---------------------------------------------------------------------------
$ws.Range("E6").Select()
